# Apply the "Updated symbol list" data refresh (prices / 1h volumes, and the
# three new coins - MXToken, LiechtensteinCryptoassetsExchange, WazirX - that
# pushed the existing GateToken/MXToken/.../LEO rows down by one slot).
#
# Every value in columns D (Price) and E (Volume(1h)) is stored as literal
# text in the source sheet (e.g. "0.9400", "317.73", "-1.11%"), so each
# assignment below is prefixed with a leading single-quote (Excel's
# "treat as text" quote-prefix) to stop COM from coercing numeric-looking
# strings into floating point Doubles and losing formatting such as
# trailing zeros. The Style reset afterwards clears the quote-prefix
# "Text" number format COM applies automatically, so cell styling stays
# untouched (matching the original, unstyled D/E cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''317.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''3.88%'
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = '''39.76'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''1.29%'
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = '''5.142'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''0.60%'
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = '''0.08226'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''1.92%'
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = '''2.049'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''5.57%'
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").Value = '''8.371'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''3.92%'
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '''0.9400'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''1.45%'
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '''0.1357'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-1.11%'
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = '''0.1993'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''3.94%'
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = '''0.09128'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''1.45%'
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = '''0.03535'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''0.40%'
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = '''0.09819'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''0.33%'
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '''0.001416'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''0.77%'
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '''0.006180'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''5.56%'
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''3.688'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''-2.14%'
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''4.315'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''2.71%'
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").Value = '''3.367'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''-0.39%'
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = '''0.3495'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''0.99%'
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").Value = '''0.1324'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''-0.09%'
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = '''4.941'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''5.72%'
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = '''0.2449'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''1.37%'
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").Value = '''0.04349'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''-0.44%'
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = '''0.001235'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''2.36%'
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = '''0.004809'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''12.39%'
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("E26").Value = '''-0.06%'
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = '''0.0003997'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''-10.13%'
$ws.Range("E27").Style = "Normal"
# Row 39
$ws.Range("D39").Value = '''0.02317'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''13.86%'
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").Value = '''0.05204'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''3.49%'
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = '''0.007745'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''2.94%'
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").Value = '''0.01046'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''7.50%'
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").Value = '''0.1412'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''5.03%'
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = '''0.002043'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-2.46%'
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = '''0.009306'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-4.87%'
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").Value = '''0.00006596'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''5.96%'
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("D47").Value = '''0.00000000749'
$ws.Range("D47").Style = "Normal"
# Row 48
$ws.Range("D48").Value = '''0.002889'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''0.55%'
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("E49").Value = '''-6.35%'
$ws.Range("E49").Style = "Normal"
